$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q2" sheet so its data survives, placing the
#    copy right after it (it will become the new "2022-Q2" sheet at index 3).
# ---------------------------------------------------------------------------
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)

# Rename the original sheet to "2022-Q4" first so the name "2022-Q2" is free,
# then rename the copy back to "2022-Q2".
$ws2.Name = "2022-Q4"
$ws3.Name = "2022-Q2"

# ---------------------------------------------------------------------------
# 2) Replace the contents of the (now) "2022-Q4" sheet with the new quarter's
#    fund-holding data.
# ---------------------------------------------------------------------------
$ws2.Cells.Clear()

$ws2.Range("B1").Value = "基金代码"
$ws2.Range("C1").Value = "基金名称"
$ws2.Range("D1").Value = "基金规模"
$ws2.Range("E1").Value = "股票总仓位"
$ws2.Range("F1").Value = "仓位占比"
$ws2.Range("G1").Value = "持有市值(亿元)"
$ws2.Range("H1").Value = "仓位排名"

$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "'004685"
$ws2.Range("C2").Value = "金元顺安元启灵活配置混合"
$ws2.Range("D2").Value = "'15.29"
$ws2.Range("E2").Value = "'76.11"
$ws2.Range("F2").Value = "'0.96"
$ws2.Range("G2").Value = "'0.1468"
$ws2.Range("H2").Value = 5

$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "'002236"
$ws2.Range("C3").Value = "大成中证360互联网+大数据100指数A"
$ws2.Range("D3").Value = "'1.15"
$ws2.Range("E3").Value = "'92.50"
$ws2.Range("F3").Value = "'0.99"
$ws2.Range("G3").Value = "'0.0114"
$ws2.Range("H3").Value = 10

$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "'003359"
$ws2.Range("C4").Value = "大成中证360互联网+大数据100指数C"
$ws2.Range("D4").Value = "'1.12"
$ws2.Range("E4").Value = "'92.50"
$ws2.Range("F4").Value = "'0.99"
$ws2.Range("G4").Value = "'0.0111"
$ws2.Range("H4").Value = 10

# The apostrophe prefix above forces these numeric-looking fund codes/ratios to
# stay text (matching the source data); ClearFormats() then drops the
# resulting quote-prefix style so the cells end up with no explicit style,
# same as every other unstyled data cell on this sheet.
$ws2.Range("B2:B4").ClearFormats()
$ws2.Range("D2:G4").ClearFormats()

# Header row + first column use the same bold/bordered style already used on
# the "总计" sheet (cellXfs index 2) - copy it over instead of re-deriving it.
$ws1.Range("B1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Update the "总计" summary sheet: push the old 2022-Q2 totals down to row
#    3 and put the new 2022-Q4 totals in row 2.
# ---------------------------------------------------------------------------
$ws1.Range("B3").Value = $ws1.Range("B2").Value2
$ws1.Range("C3").Value = $ws1.Range("C2").Value2
$ws1.Range("D3").Value = $ws1.Range("D2").Value2

$ws1.Range("A2").Copy()
$ws1.Range("A3").PasteSpecial(-4122)
$ws1.Range("A3").Value = 1

$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 3
$ws1.Range("D2").Value = 0.17
